# Implement login page ui elements
# Appends three new data rows (106-108) to each of the four worksheets,
# mirroring the existing daily-log row pattern (columns A-I).

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# Per-sheet, per-row data extracted from the target change.
# Each row entry: A (serial date), B, C, D, E (text/hex strings),
# F, G, H, I (numeric values).
$sheetsData = @(
    @{
        Name = "MID_LFT_#1"
        Rows = @(
            @{ Row = 106; A = 45892.46393518519; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; D = "0x01,0x08"; E = "0x07"; F = 400; G = "5.68631262647113e+23"; H = 264; I = 7 }
            @{ Row = 107; A = 45893.46261574074; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; D = "0x01,0x08"; E = "0x07"; F = 400; G = "5.68631262647113e+23"; H = 264; I = 7 }
            @{ Row = 108; A = 45894.46487268519; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; D = "0x01,0x08"; E = "0x07"; F = 400; G = "5.68631262647113e+23"; H = 264; I = 7 }
        )
    }
    @{
        Name = "MID_LFT_#2"
        Rows = @(
            @{ Row = 106; A = 45892.46393518519; B = "0x01,0x7c"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x18"; E = "0x19"; F = 380; G = "5.68432987514711e+23"; H = 280; I = 25 }
            @{ Row = 107; A = 45893.46261574074; B = "0x01,0x7c"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x14"; E = "0x19"; F = 380; G = "5.68432987514711e+23"; H = 276; I = 25 }
            @{ Row = 108; A = 45894.46487268519; B = "0x01,0x7c"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x14"; E = "0x19"; F = 380; G = "5.68432987514711e+23"; H = 276; I = 25 }
        )
    }
    @{
        Name = "MID_PLT_#1"
        Rows = @(
            @{ Row = 106; A = 45892.46393518519; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x00,0x5B"; E = "0x15"; F = 110; G = "5.68631262647113e+23"; H = 91; I = 15 }
            @{ Row = 107; A = 45893.46261574074; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x00,0x5B"; E = "0x15"; F = 110; G = "5.68631262647113e+23"; H = 91; I = 15 }
            @{ Row = 108; A = 45894.46487268519; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x00,0x5B"; E = "0x15"; F = 110; G = "5.68631262647113e+23"; H = 91; I = 15 }
        )
    }
    @{
        Name = "MID_PLT_#2"
        Rows = @(
            @{ Row = 106; A = 45892.46393518519; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x00,0x72"; E = "0x9"; F = 130; G = "5.68631262647113e+23"; H = 114; I = 9 }
            @{ Row = 107; A = 45893.46261574074; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x00,0x72"; E = "0x9"; F = 130; G = "5.68631262647113e+23"; H = 114; I = 9 }
            @{ Row = 108; A = 45894.46487268519; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x00,0x72"; E = "0x9"; F = 130; G = "5.68631262647113e+23"; H = 114; I = 9 }
        )
    }
)

foreach ($sheetData in $sheetsData) {
    $ws = $wb.Worksheets.Item($sheetData.Name)

    foreach ($rowData in $sheetData.Rows) {
        $r = $rowData.Row

        $ws.Cells.Item($r, 1).Value = $rowData.A
        $ws.Cells.Item($r, 1).NumberFormat = $dateFormat

        $ws.Cells.Item($r, 2).Value = $rowData.B
        $ws.Cells.Item($r, 3).Value = $rowData.C
        $ws.Cells.Item($r, 4).Value = $rowData.D
        $ws.Cells.Item($r, 5).Value = $rowData.E

        $ws.Cells.Item($r, 6).Value = $rowData.F
        $ws.Cells.Item($r, 7).Value = [double]$rowData.G
        $ws.Cells.Item($r, 8).Value = $rowData.H
        $ws.Cells.Item($r, 9).Value = $rowData.I
    }
}
